$d = $word.ActiveDocument

# 1. Mark the "[show the other testing class...]" direction as struck-through,
#    matching the other completed-direction annotations in the document.
$rng1 = $d.Content
$found1 = $rng1.Find.Execute(
    "[show the other testing class and each of the methods within the class being used]",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found1) {
    $rng1.Font.StrikeThrough = 1
}

# 2. Move the "_GoBack" last-edit bookmark from the end of the "[show log file]"
#    paragraph to mid-word inside "...whenever I need to prompt the user...",
#    splitting "prompt" into "p" + "rompt" around the new bookmark location.
$rng2 = $d.Content
$found2 = $rng2.Find.Execute(
    "have used error handling in my code whenever I need to p",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found2) {
    $insertionPoint = $rng2.Duplicate
    $insertionPoint.Collapse(0)
    $d.Bookmarks.Add("_GoBack", $insertionPoint)
}
